{"js": "// Update a handful of statistics in the \"Republican votes region\" table\n// (West and Midwest rows) to reflect corrected figures from the final\n// pass of the descriptive-statistics table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row/col are 0-based indices into the table grid (row 0 is the header row).\nconst changes = [\n  { row: 2, col: 2, oldText: \"0.5756993\", newText: \"0.5756981\" },   // West.mean\n  { row: 2, col: 3, oldText: \"0.1720934\", newText: \"0.1720930\" },   // West.sd\n  { row: 2, col: 5, oldText: \"0.5857550\", newText: \"0.5857534\" },   // West.trimmed\n  { row: 2, col: 10, oldText: \"-0.4549777\", newText: \"-0.4549619\" }, // West.skew\n  { row: 2, col: 11, oldText: \"-0.5963621\", newText: \"-0.5963517\" }, // West.kurtosis\n  { row: 2, col: 12, oldText: \"0.008387320\", newText: \"0.008387298\" }, // West.se\n  { row: 4, col: 1, oldText: \"1,053\", newText: \"1,054\" },           // Midwest.n\n  { row: 4, col: 2, oldText: \"0.6621178\", newText: \"0.6615655\" },   // Midwest.mean\n  { row: 4, col: 3, oldText: \"0.1224697\", newText: \"0.1237178\" },   // Midwest.sd\n  { row: 4, col: 5, oldText: \"0.6703084\", newText: \"0.6701066\" },   // Midwest.trimmed\n  { row: 4, col: 7, oldText: \"0.16\", newText: \"0.08\" },             // Midwest.min\n  { row: 4, col: 9, oldText: \"0.77\", newText: \"0.85\" },             // Midwest.range\n  { row: 4, col: 10, oldText: \"-0.7119323\", newText: \"-0.7753996\" }, // Midwest.skew\n  { row: 4, col: 11, oldText: \"0.7338338\", newText: \"1.0330871\" },  // Midwest.kurtosis\n  { row: 4, col: 12, oldText: \"0.003774111\", newText: \"0.003810763\" }, // Midwest.se\n];\n\nfor (const change of changes) {\n  const cell = table.getCell(change.row, change.col);\n  cell.body.load(\"text\");\n  await context.sync();\n\n  const actual = cell.body.text.trim();\n  if (actual !== change.oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${change.row}, col ${change.col}: ` +\n      `expected \"${change.oldText}\" but found \"${actual}\"`\n    );\n  }\n\n  // Replace only the text of the cell's range so paragraph/run formatting\n  // (font, size, color, alignment, spacing, indent) is preserved.\n  const range = cell.body.getRange();\n  range.insertText(change.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update a handful of statistics in the \"Republican votes region\" table\n# (West and Midwest rows) to reflect corrected figures from the final\n# pass of the descriptive-statistics table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row/Col are 1-based table coordinates (row 1 is the header row).\n$changes = @(\n  @{ Row = 3; Col = 3;  OldText = '0.5756993';   NewText = '0.5756981' },   # West.mean\n  @{ Row = 3; Col = 4;  OldText = '0.1720934';   NewText = '0.1720930' },   # West.sd\n  @{ Row = 3; Col = 6;  OldText = '0.5857550';   NewText = '0.5857534' },   # West.trimmed\n  @{ Row = 3; Col = 11; OldText = '-0.4549777';  NewText = '-0.4549619' },  # West.skew\n  @{ Row = 3; Col = 12; OldText = '-0.5963621';  NewText = '-0.5963517' },  # West.kurtosis\n  @{ Row = 3; Col = 13; OldText = '0.008387320'; NewText = '0.008387298' }, # West.se\n  @{ Row = 5; Col = 2;  OldText = '1,053';       NewText = '1,054' },       # Midwest.n\n  @{ Row = 5; Col = 3;  OldText = '0.6621178';   NewText = '0.6615655' },   # Midwest.mean\n  @{ Row = 5; Col = 4;  OldText = '0.1224697';   NewText = '0.1237178' },   # Midwest.sd\n  @{ Row = 5; Col = 6;  OldText = '0.6703084';   NewText = '0.6701066' },   # Midwest.trimmed\n  @{ Row = 5; Col = 8;  OldText = '0.16';        NewText = '0.08' },        # Midwest.min\n  @{ Row = 5; Col = 10; OldText = '0.77';        NewText = '0.85' },        # Midwest.range\n  @{ Row = 5; Col = 11; OldText = '-0.7119323';  NewText = '-0.7753996' },  # Midwest.skew\n  @{ Row = 5; Col = 12; OldText = '0.7338338';   NewText = '1.0330871' },   # Midwest.kurtosis\n  @{ Row = 5; Col = 13; OldText = '0.003774111'; NewText = '0.003810763' }  # Midwest.se\n)\n\nforeach ($change in $changes) {\n  $cell = $t.Cell($change.Row, $change.Col)\n  $r = $cell.Range\n  # Drop the trailing end-of-cell mark so we only touch the visible text,\n  # which lets Word reuse the existing run (and its formatting).\n  $r.MoveEnd(1, -1) | Out-Null\n\n  if ($r.Text -ne $change.OldText) {\n    throw \"Unexpected cell text at row $($change.Row), col $($change.Col): expected '$($change.OldText)' but found '$($r.Text)'\"\n  }\n\n  $r.Text = $change.NewText\n}\n"}
